# SRS first review - Answer SIQ opened questions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update the Answer in row 8 (I8) with the refined wording agreed on review
$ws.Range("I8").Value = "it should not allow it if it is multiplication or division od addition but should allow it if it is just subtraction as it should mean a sign in such case."

# Fill in the previously-empty Return Date for row 10 (H10), using the same
# date format as the other Return Date cells in this column (H7:H9)
$ws.Range("H10").NumberFormat = "m/d/yy"
$ws.Range("H10").Value = 44014

# Scroll the view and move the active selection to where the reviewer left off
$win = $wb.Windows.Item(1)
$win.ScrollColumn = 3
$win.ScrollRow = 7
$ws.Range("H10").Select()
